$wb = $excel.ActiveWorkbook

# ALC row 74
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 3415
$ws.Range("I74").Value = 2624.375
$ws.Range("J74").Value = 4680
$ws.Range("K74").Value = 2624.375
$ws.Range("L74").Value = 4680
$ws.Range("M74").Value = -1688.375
$ws.Range("N74").Value = -6552

# ALC row 77
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H77").Value = 3415
$ws.Range("I77").Value = 2624.375
$ws.Range("J77").Value = 4680
$ws.Range("K77").Value = 13121.875
$ws.Range("L77").Value = 23400
$ws.Range("M77").Value = -8441.875
$ws.Range("N77").Value = -32760

# ALC row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 2940.9688
$ws.Range("I132").Value = 3161.1072
$ws.Range("K132").Value = 9483.321599999999
$ws.Range("M132").Value = -6953.321599999999

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2498.3044
$ws.Range("I137").Value = 2697.2144
$ws.Range("J137").Value = 2188.889
$ws.Range("K137").Value = 8091.6432
$ws.Range("L137").Value = 6566.667
$ws.Range("M137").Value = -5541.6432
$ws.Range("N137").Value = -11666.667

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 13701354
$ws.Range("I138").Value = 35715588
$ws.Range("J138").Value = 3608.889
$ws.Range("K138").Value = 107146764
$ws.Range("L138").Value = 10826.667
$ws.Range("M138").Value = -107141624
$ws.Range("N138").Value = -21106.667

# ALC row 141
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 1348.738
$ws.Range("I141").Value = 745.3823
$ws.Range("J141").Value = 3913
$ws.Range("K141").Value = 2236.1469
$ws.Range("L141").Value = 11739
$ws.Range("M141").Value = 2943.8531
$ws.Range("N141").Value = -22099

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4252.459
$ws.Range("I32").Value = 2844
$ws.Range("K32").Value = 2844
$ws.Range("M32").Value = -2557

# ARM row 49
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H49").Value = 12000
$ws.Range("J49").Value = 12000
$ws.Range("L49").Value = 12000
$ws.Range("N49").Value = -12520

# ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 392615.9
$ws.Range("I61").Value = 440077.4
$ws.Range("J61").Value = 3431.6
$ws.Range("K61").Value = 440077.4
$ws.Range("L61").Value = 3431.6
$ws.Range("M61").Value = -439865.4
$ws.Range("N61").Value = -3855.6

# ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 21278186
$ws.Range("I74").Value = 24391792
$ws.Range("J74").Value = 1883
$ws.Range("K74").Value = 24391792
$ws.Range("L74").Value = 1883
$ws.Range("M74").Value = -24390918
$ws.Range("N74").Value = -3631

# ARM row 76
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H76").Value = 19998
$ws.Range("J76").Value = 19998
$ws.Range("L76").Value = 19998
$ws.Range("N76").Value = -20674

# ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 21278186
$ws.Range("I77").Value = 24391792
$ws.Range("J77").Value = 1883
$ws.Range("K77").Value = 121958960
$ws.Range("L77").Value = 9415
$ws.Range("M77").Value = -121954592
$ws.Range("N77").Value = -18151

# ARM row 79
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H79").Value = 19998
$ws.Range("J79").Value = 19998
$ws.Range("L79").Value = 19998
$ws.Range("N79").Value = -22338

# ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 14885.895
$ws.Range("I132").Value = 1801.931
$ws.Range("J132").Value = 57045.332
$ws.Range("K132").Value = 5405.793
$ws.Range("L132").Value = 171135.996
$ws.Range("M132").Value = -2875.793
$ws.Range("N132").Value = -176195.996

# ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 392615.9
$ws.Range("I136").Value = 440077.4
$ws.Range("J136").Value = 3431.6
$ws.Range("K136").Value = 1320232.2
$ws.Range("L136").Value = 10294.8
$ws.Range("M136").Value = -1317682.2
$ws.Range("N136").Value = -15394.8

# BSM row 99
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1885
$ws.Range("I99").Value = 1952.5
$ws.Range("J99").Value = 1750
$ws.Range("K99").Value = 1952.5
$ws.Range("L99").Value = 1750
$ws.Range("M99").Value = -454.5
$ws.Range("N99").Value = -4746

# BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2327.5173
$ws.Range("I134").Value = 2389.8125
$ws.Range("J134").Value = 2028.5
$ws.Range("K134").Value = 7169.4375
$ws.Range("L134").Value = 6085.5
$ws.Range("M134").Value = -4634.4375
$ws.Range("N134").Value = -11155.5

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3725.818
$ws.Range("I31").Value = 2215.35
$ws.Range("J31").Value = 6049.615
$ws.Range("K31").Value = 2215.35
$ws.Range("L31").Value = 6049.615
$ws.Range("M31").Value = -1920.35
$ws.Range("N31").Value = -6639.615

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 3725.818
$ws.Range("I34").Value = 2215.35
$ws.Range("J34").Value = 6049.615
$ws.Range("K34").Value = 2215.35
$ws.Range("L34").Value = 6049.615
$ws.Range("M34").Value = -2013.35
$ws.Range("N34").Value = -6453.615

# CRP row 52
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H52").Value = 40000
$ws.Range("J52").Value = 40000
$ws.Range("L52").Value = 40000
$ws.Range("N52").Value = -40588

# CRP row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 7023.375
$ws.Range("I58").Value = 628.92725
$ws.Range("J58").Value = 21091.16
$ws.Range("K58").Value = 628.92725
$ws.Range("L58").Value = 21091.16
$ws.Range("M58").Value = -425.92725
$ws.Range("N58").Value = -21497.16

# CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1479.8387
$ws.Range("I132").Value = 1209.1852
$ws.Range("J132").Value = 3306.75
$ws.Range("K132").Value = 3627.5556
$ws.Range("L132").Value = 9920.25
$ws.Range("M132").Value = -1097.5556
$ws.Range("N132").Value = -14980.25

# CRP row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 696.9048
$ws.Range("I134").Value = 700
$ws.Range("J134").Value = 688.1818
$ws.Range("K134").Value = 2100
$ws.Range("L134").Value = 2064.5454
$ws.Range("M134").Value = 435
$ws.Range("N134").Value = -7134.5454

# CRP row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 7023.375
$ws.Range("I136").Value = 628.92725
$ws.Range("J136").Value = 21091.16
$ws.Range("K136").Value = 1886.78175
$ws.Range("L136").Value = 63273.48
$ws.Range("M136").Value = 663.2182500000001
$ws.Range("N136").Value = -68373.48

# CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 771.59
$ws.Range("J131").Value = 787.54254
$ws.Range("L131").Value = 2362.62762
$ws.Range("N131").Value = -12442.62762

# CUL row 137
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 19610754
$ws.Range("J137").Value = 27781482
$ws.Range("L137").Value = 83344446
$ws.Range("N137").Value = -83354646

# GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 24055.084
$ws.Range("I132").Value = 3505.8096
$ws.Range("J132").Value = 167900
$ws.Range("K132").Value = 10517.4288
$ws.Range("L132").Value = 503700
$ws.Range("M132").Value = -7987.4288
$ws.Range("N132").Value = -508760

# GSM row 136
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H136").Value = 14313
$ws.Range("J136").Value = 14313
$ws.Range("L136").Value = 42939
$ws.Range("N136").Value = -48039

# GSM row 137
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H137").Value = 60000
$ws.Range("J137").Value = 60000
$ws.Range("L137").Value = 60000
$ws.Range("N137").Value = -70200

# LTW row 7
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5690.65
$ws.Range("I7").Value = 4263.5
$ws.Range("J7").Value = 6642.0835
$ws.Range("K7").Value = 4263.5
$ws.Range("L7").Value = 6642.0835
$ws.Range("M7").Value = -4151.5
$ws.Range("N7").Value = -6866.0835

# LTW row 22
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2084.6924
$ws.Range("J22").Value = 3050
$ws.Range("L22").Value = 3050
$ws.Range("N22").Value = -3640

# LTW row 27
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 2084.6924
$ws.Range("J27").Value = 3050
$ws.Range("L27").Value = 3050
$ws.Range("N27").Value = -3264

# LTW row 40
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3732.9565
$ws.Range("I40").Value = 3619.7222
$ws.Range("J40").Value = 4140.6
$ws.Range("K40").Value = 3619.7222
$ws.Range("L40").Value = 4140.6
$ws.Range("M40").Value = -3483.7222
$ws.Range("N40").Value = -4412.6

# LTW row 41
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H41").Value = 17000
$ws.Range("J41").Value = 17000
$ws.Range("L41").Value = 17000
$ws.Range("N41").Value = -17876

# LTW row 46
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1499.2858
$ws.Range("I46").Value = 999
$ws.Range("K46").Value = 999
$ws.Range("M46").Value = -811

# LTW row 47
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H47").Value = 29990
$ws.Range("J47").Value = 29990
$ws.Range("L47").Value = 29990
$ws.Range("N47").Value = -30970

# LTW row 52
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H52").Value = 29990
$ws.Range("J52").Value = 29990
$ws.Range("L52").Value = 29990
$ws.Range("N52").Value = -30456

# LTW row 68
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2214.8235
$ws.Range("I68").Value = 2172.9092
$ws.Range("J68").Value = 2291.6667
$ws.Range("K68").Value = 2172.9092
$ws.Range("L68").Value = 2291.6667
$ws.Range("M68").Value = -1423.9092
$ws.Range("N68").Value = -3789.6667

# LTW row 71
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 2214.8235
$ws.Range("I71").Value = 2172.9092
$ws.Range("J71").Value = 2291.6667
$ws.Range("K71").Value = 10864.546
$ws.Range("L71").Value = 11458.3335
$ws.Range("M71").Value = -7120.546
$ws.Range("N71").Value = -18946.3335

# LTW row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 855469.4399999999
$ws.Range("I122").Value = 1963649.8
$ws.Range("J122").Value = 3023
$ws.Range("K122").Value = 5890949.4
$ws.Range("L122").Value = 9069
$ws.Range("M122").Value = -5888499.4
$ws.Range("N122").Value = -13969

# LTW row 126
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 5690.65
$ws.Range("I126").Value = 4263.5
$ws.Range("J126").Value = 6642.0835
$ws.Range("K126").Value = 12790.5
$ws.Range("L126").Value = 19926.2505
$ws.Range("M126").Value = -10320.5
$ws.Range("N126").Value = -24866.2505

# LTW row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 1062.5254
$ws.Range("I132").Value = 1064.4259
$ws.Range("J132").Value = 1042
$ws.Range("K132").Value = 3193.2777
$ws.Range("L132").Value = 3126
$ws.Range("M132").Value = -663.2776999999996
$ws.Range("N132").Value = -8186

# LTW row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 943.8570999999999
$ws.Range("I136").Value = 948
$ws.Range("J136").Value = 899.6667
$ws.Range("K136").Value = 2844
$ws.Range("L136").Value = 2699.0001
$ws.Range("M136").Value = -294
$ws.Range("N136").Value = -7799.0001

# WVR row 62
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4001
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0

# WVR row 65
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 4001
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0

# WVR row 107
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 3247347.5
$ws.Range("I107").Value = 610.7273
$ws.Range("K107").Value = 1832.1819
$ws.Range("M107").Value = 87.81809999999996

# WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 578.90625
$ws.Range("I132").Value = 465.1607
$ws.Range("J132").Value = 1375.125
$ws.Range("K132").Value = 1395.4821
$ws.Range("L132").Value = 4125.375
$ws.Range("M132").Value = 1134.5179
$ws.Range("N132").Value = -9185.375

# WVR row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 12988637
$ws.Range("I136").Value = 19608650
$ws.Range("J136").Value = 3229.423
$ws.Range("K136").Value = 58825950
$ws.Range("L136").Value = 9688.269
$ws.Range("M136").Value = -58823400
$ws.Range("N136").Value = -14788.269
